# Bill of Materials - Added radiation shield parts
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Instrument")

# --- Row 15: Radiation Shield ---
$ws.Range("B15").Value = "Radiation Shield"
$ws.Range("C15").Value = "N/A"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 6
$ws.Range("G15").Value = 6
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = "UofU MechE Dept."
$ws.Range("L15").Value = "N/A"
$ws.Range("M15").Value = "N/A"
$ws.Range("O15").Value = "Printed from UofU MechE 3D printer. Maybe cheaper elsewhere?"

# --- Row 16: Set Screw ---
$ws.Range("B16").Value = "Set Screw"
$ws.Range("C16").Value = "N/A"

# header label change: Quantity/Board -> Quantity/Instrument
$ws.Range("D3").Value = "Quantity/Instrument"

$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0.0778
$ws.Range("G16").Value = 0.0778
$ws.Range("I16").Value = 0.0778
$ws.Range("K16").Value = "McMaster-Carr"
$ws.Range("L16").Value = "91290A111"
$ws.Range("M16").Value = "M3 0.5 x 6 Socket Head"
$ws.Range("N16").Value = "http://www.mcmaster.com/#91290A111"
$ws.Range("O16").Value = "Comes in packs of 100 @ `$7.78 per pack"

# --- Row 17: Plasti-Dip White Spray ---
$ws.Range("B17").Value = "Plasti-Dip White Spray"
$ws.Range("C17").Value = "N/A"
$ws.Range("D17").Formula = "=1/200"
$ws.Range("E17").Value = 5.98
$ws.Range("G17").Value = 5.98
$ws.Range("I17").Value = 5.98
$ws.Range("K17").Value = "Home Depot"
$ws.Range("L17").Value = 203286992
$ws.Range("M17").Value = "11 oz. Can"
$ws.Range("N17").Value = "http://www.homedepot.com/p/Plasti-Dip-11-oz-White-General-Purpose-Rubber-Coating-Spray-11207-6/203286992"
$ws.Range("O17").Value = "Assume 200 BTEMS can be made with one can"

# L17 keeps its numeric value but switches the cell format to Text (matches the
# authored workbook: style index with numFmtId 49 "@")
$ws.Range("L17").NumberFormat = "@"

# --- Cosmetic: column width + style tweaks on L7 / L10 to match authoring session ---
$ws.Columns.Item(4).ColumnWidth = 17.1
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L10").NumberFormat = "@"

# --- selection / window state ---
$ws.Activate()
$ws.Range("C34").Select()
